$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52, shifting existing rows 52:128 down to 53:129.
$ws.Rows(52).Insert()

# Populate the newly inserted row 52 with the new data record.
$ws.Range("A52").Value = 9
$ws.Range("B52").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C52").Value = "Metropolitana"
$ws.Range("D52").Value = 44792
$ws.Range("E52").Value = 13
$ws.Range("F52").Value = 100112022
$ws.Range("G52").Value = "Arveja Verde"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 36
$ws.Range("K52").Value = 40000
$ws.Range("L52").Value = 42000
$ws.Range("M52").Value = 41000
$ws.Range("N52").Value = "$/malla 25 kilos"
$ws.Range("O52").Value = "Provincia de Huasco"
$ws.Range("P52").Value = 1640
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"
